# Update cryptos list: apply per-cell text values, coercing to text
# (NumberFormat "@" + ClearFormats) so numeric-looking strings stay text
# cells exactly like the source inlineStr cells, without leaving a stray style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "20.514.34"
Set-TextValue "E2" "  +2.45%  "

Set-TextValue "D3" "1.471.34"
Set-TextValue "E3" "  +3.51%  "

Set-TextValue "D4" "1.008"
Set-TextValue "E4" "  +0.79%  "

Set-TextValue "D5" "0.9435"
Set-TextValue "E5" "  -5.67%  "

Set-TextValue "D6" "280.94"
Set-TextValue "E6" "  +2.51%  "

Set-TextValue "D7" "0.3714"
Set-TextValue "E7" "  -0.03%  "

Set-TextValue "D8" "0.3194"
Set-TextValue "E8" "  +3.86%  "

Set-TextValue "D9" "41.40"
Set-TextValue "E9" "  +4.33%  "

Set-TextValue "E10" "  +4.70%  "

Set-TextValue "D11" "0.06673"
Set-TextValue "E11" "  +1.04%  "

Set-TextValue "E12" "  +0.28%  "

Set-TextValue "D13" "5.597"
Set-TextValue "E13" "  +3.28%  "

Set-TextValue "D14" "18.28"
Set-TextValue "E14" "  +6.67%  "

Set-TextValue "D15" "6.235"
Set-TextValue "E15" "  +0.98%  "

Set-TextValue "D16" "1.479.44"
Set-TextValue "E16" "  +4.10%  "

Set-TextValue "D17" "0.00001034"
Set-TextValue "E17" "  +2.56%  "

Set-TextValue "D18" "0.9435"
Set-TextValue "E18" "  -5.68%  "

Set-TextValue "D19" "0.05733"
Set-TextValue "E19" "  -1.44%  "

Set-TextValue "D20" "72.16"
Set-TextValue "E20" "  -3.36%  "

Set-TextValue "D21" "5.694"
Set-TextValue "E21" "  +0.82%  "

Set-TextValue "D22" "14.77"

Set-TextValue "D23" "11.24"
Set-TextValue "E23" "  +1.71%  "

Set-TextValue "D24" "2.273"
Set-TextValue "E24" "  -2.58%  "

Set-TextValue "D25" "20.796.22"
Set-TextValue "E25" "  +3.80%  "

Set-TextValue "D26" "2.304"
Set-TextValue "E26" "  +0.21%  "

Set-TextValue "D27" "138.21"
Set-TextValue "E27" "  -0.44%  "

Set-TextValue "D28" "17.60"
Set-TextValue "E28" "  +4.22%  "

Set-TextValue "D29" "1.638.81"
Set-TextValue "E29" "  +3.66%  "

Set-TextValue "D30" "113.70"
Set-TextValue "E30" "  +4.23%  "

Set-TextValue "D31" "3.936"
Set-TextValue "E31" "  +3.15%  "

Set-TextValue "E32" "  -1.99%  "

Set-TextValue "D33" "0.8517"
Set-TextValue "E33" "  -4.15%  "

Set-TextValue "D34" "1.612"
Set-TextValue "E34" "  +26.75%  "

Set-TextValue "D35" "0.07834"
Set-TextValue "E35" "  +1.18%  "

Set-TextValue "D36" "0.06058"
Set-TextValue "E36" "  +5.96%  "

Set-TextValue "D37" "4.944"
Set-TextValue "E37" "  +3.25%  "

Set-TextValue "D38" "10.71"
Set-TextValue "E38" "  -5.37%  "

Set-TextValue "D39" "0.02073"
Set-TextValue "E39" "  +1.84%  "

Set-TextValue "D40" "1.123"
Set-TextValue "E40" "  +2.88%  "

Set-TextValue "D41" "0.1907"

Set-TextValue "D42" "0.9577"
Set-TextValue "E42" "  -4.25%  "

Set-TextValue "D43" "7.562"
Set-TextValue "E43" "  -10.38%  "

Set-TextValue "D44" "0.5418"
Set-TextValue "E44" "  +1.71%  "

Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "12.51"
Set-TextValue "E45" "  +1.66%  "

Set-TextValue "B46" "PancakeSwap"
Set-TextValue "C46" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D46" "3.590"
Set-TextValue "E46" "  +1.48%  "

Set-TextValue "D47" "122.38"
Set-TextValue "E47" "  +11.59%  "

Set-TextValue "D48" "0.5324"
Set-TextValue "E48" "  +3.66%  "

Set-TextValue "D49" "1.829"
Set-TextValue "E49" "  +1.51%  "

Set-TextValue "D50" "0.06458"
Set-TextValue "E50" "  +4.33%  "

Set-TextValue "D51" "1.047"
Set-TextValue "E51" "  -0.21%  "

